$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.665.15'
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").Value = '1.893.09'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = '311.68'
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").Value = '0.4908'
$ws.Range("E7").Value = '  +0.63%  '
$ws.Range("D8").Value = '0.3794'
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07310'
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("E10").Value = '  -5.16%  '
$ws.Range("D11").Value = '20.54'
$ws.Range("E11").Value = '  -3.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07650'
$ws.Range("E12").Value = '  -2.39%  '
$ws.Range("D13").Value = '1.915.70'
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("D14").Value = '5.461'
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("D15").Value = '6.631'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").Value = '90.97'
$ws.Range("E16").Value = '  -1.66%  '
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008752'
$ws.Range("E18").Value = '  -2.12%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").Value = '27.785.01'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").Value = '14.44'
$ws.Range("E21").Value = '  -4.43%  '
$ws.Range("D23").Value = '2.152.07'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '10.75'
$ws.Range("E24").Value = '  -2.22%  '
$ws.Range("D25").Value = '154.05'
$ws.Range("E25").Value = '  -1.95%  '
$ws.Range("D26").Value = '1.854'
$ws.Range("E26").Value = '  -5.42%  '
$ws.Range("D27").Value = '18.37'
$ws.Range("E27").Value = '  -2.06%  '
$ws.Range("D28").Value = '2.156'
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("D29").Value = '115.18'
$ws.Range("E29").Value = '  -1.54%  '
$ws.Range("E30").Value = '  -3.81%  '
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("D32").Value = '3.207'
$ws.Range("E32").Value = '  -3.80%  '
$ws.Range("D33").Value = '1.224'
$ws.Range("E33").Value = '  -2.50%  '
$ws.Range("D34").Value = '0.7618'
$ws.Range("E34").Value = '  -3.06%  '
$ws.Range("D35").Value = '4.621'
$ws.Range("E35").Value = '  -2.16%  '
$ws.Range("D36").Value = '2.556'
$ws.Range("E36").Value = '  -8.71%  '
$ws.Range("D37").Value = '0.02034'
$ws.Range("E37").Value = '  -1.46%  '
$ws.Range("D38").Value = '1.098'
$ws.Range("E38").Value = '  -3.45%  '
$ws.Range("D39").Value = '0.05291'
$ws.Range("E39").Value = '  -2.33%  '
$ws.Range("D40").Value = '2.986'
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("D41").Value = '0.5455'
$ws.Range("E41").Value = '  -2.96%  '
$ws.Range("D42").Value = '6.874'
$ws.Range("E42").Value = '  -3.26%  '
$ws.Range("D45").Value = '0.1519'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").Value = '10.57'
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("D47").Value = '0.4783'
$ws.Range("E47").Value = '  -3.90%  '
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("D49").Value = '1.628'
$ws.Range("E49").Value = '  -3.81%  '
$ws.Range("D50").Value = '67.31'
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("E51").Value = '  -1.20%  '

# Rows 43/44: Quant and Aptos swap positions, with updated price/volume data
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '8.519'
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '112.79'
$ws.Range("E44").Value = '  +5.57%  '
